$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUMIFS-VLOOKUP")

# Update the worksheet header text.
$ws.Range("A1").Value = "Excel Skills demonstrated:"

# Rename the client codes used throughout the client-list / data-entry
# table to the new generic "Client A".."Client E" labels.
# Mapping: TMMG -> Client A, LOC -> Client B, Greatland -> Client C,
#          Global -> Client D, Kaiser -> Client E
$ws.Range("B21").Value = "Client A"
$ws.Range("G21").Value = "Client A"

$ws.Range("B22").Value = "Client B"
$ws.Range("G22").Value = "Client B"

$ws.Range("B23").Value = "Client C"
$ws.Range("G23").Value = "Client C"

$ws.Range("B24").Value = "Client A"
$ws.Range("G24").Value = "Client D"

$ws.Range("B25").Value = "Client D"
$ws.Range("G25").Value = "Client E"

$ws.Range("B26").Value = "Client C"
$ws.Range("B27").Value = "Client D"
$ws.Range("B28").Value = "Client A"

$ws.Range("B29").Value = "Client D"
$ws.Range("B29").Style = "Normal"

$ws.Range("B30").Value = "Client A"

# Remove the now unused external workbook link (Financial/2021 Income.xlsx).
$wb.BreakLink("Financial/2021%20Income.xlsx", 1)

# Update the active selection on the worksheet.
$ws.Activate()
$ws.Range("G34").Select()
